# Regenerate save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals. This updates column G ("K") values for rows 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 1
    4  = 2
    5  = 3
    6  = 6
    7  = 6
    8  = 2
    9  = 8
    10 = 5
    11 = 3
    12 = 7
    13 = 3
    14 = 3
    15 = 2
    16 = 4
    17 = 8
    18 = 5
    19 = 6
    20 = 6
    21 = 3
    22 = 0
    23 = 4
    24 = 9
    25 = 1
    26 = 2
    27 = 6
    28 = 2
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
